$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '29.811.44'
Set-TextValue $ws 'E2' '  -0.47%  '

# Row 3
Set-TextValue $ws 'D3' '1.867.89'
Set-TextValue $ws 'E3' '  -1.34%  '

# Row 4
Set-TextValue $ws 'D4' '1.000'
Set-TextValue $ws 'E4' '  -0.05%  '

# Row 5
Set-TextValue $ws 'D5' '0.7321'
Set-TextValue $ws 'E5' '  -5.01%  '

# Row 6
Set-TextValue $ws 'D6' '241.73'
Set-TextValue $ws 'E6' '  -1.13%  '

# Row 7
Set-TextValue $ws 'D7' '0.9997'
Set-TextValue $ws 'E7' '  -0.09%  '

# Row 8
Set-TextValue $ws 'D8' '0.3137'
Set-TextValue $ws 'E8' '  +0.29%  '

# Row 9
Set-TextValue $ws 'B9' 'Solana'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws 'D9' '24.57'
Set-TextValue $ws 'E9' '  -4.36%  '

# Row 10
Set-TextValue $ws 'B10' 'Dogecoin'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws 'D10' '0.07066'
Set-TextValue $ws 'E10' '  -2.81%  '

# Row 11
Set-TextValue $ws 'D11' '0.08415'
Set-TextValue $ws 'E11' '  +4.69%  '

# Row 12
Set-TextValue $ws 'B12' 'Polygon'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D12' '0.7473'
Set-TextValue $ws 'E12' '  -3.18%  '

# Row 13
Set-TextValue $ws 'B13' 'Polkadot'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D13' '5.347'
Set-TextValue $ws 'E13' '  -1.88%  '

# Row 14
Set-TextValue $ws 'B14' 'WrappedEther'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D14' '1.844.53'
Set-TextValue $ws 'E14' '  -2.48%  '

# Row 15
Set-TextValue $ws 'D15' '91.95'
Set-TextValue $ws 'E15' '  -3.25%  '

# Row 16
Set-TextValue $ws 'D16' '29.813.84'
Set-TextValue $ws 'E16' '  -0.25%  '

# Row 17
Set-TextValue $ws 'D17' '6.031'
Set-TextValue $ws 'E17' '  -2.81%  '

# Row 18
Set-TextValue $ws 'D18' '13.50'
Set-TextValue $ws 'E18' '  -3.41%  '

# Row 19
Set-TextValue $ws 'D19' '240.66'
Set-TextValue $ws 'E19' '  -2.33%  '

# Row 20
Set-TextValue $ws 'D20' '0.000007787'
Set-TextValue $ws 'E20' '  -0.77%  '

# Row 21
Set-TextValue $ws 'D21' '0.9999'
Set-TextValue $ws 'E21' '  -0.05%  '

# Row 22
Set-TextValue $ws 'D22' '2.129.59'
Set-TextValue $ws 'E22' '  +0.52%  '

# Row 23
Set-TextValue $ws 'B23' 'Chainlink'
Set-TextValue $ws 'C23' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D23' '7.917'
Set-TextValue $ws 'E23' '  -2.24%  '

# Row 24
Set-TextValue $ws 'B24' 'BinanceUSD'
Set-TextValue $ws 'C24' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws 'D24' '1.001'
Set-TextValue $ws 'E24' '  -0.01%  '

# Row 25
Set-TextValue $ws 'D25' '0.1561'
Set-TextValue $ws 'E25' '  -1.43%  '

# Row 26
Set-TextValue $ws 'D26' '9.292'
Set-TextValue $ws 'E26' '  -2.27%  '

# Row 27
Set-TextValue $ws 'D27' '163.60'
Set-TextValue $ws 'E27' '  +0.85%  '

# Row 28
Set-TextValue $ws 'D28' '18.58'
Set-TextValue $ws 'E28' '  -1.01%  '

# Row 29
Set-TextValue $ws 'D29' '2.017'
Set-TextValue $ws 'E29' '  -0.85%  '

# Row 30
Set-TextValue $ws 'D30' '1.465'
Set-TextValue $ws 'E30' '  +3.63%  '

# Row 31
Set-TextValue $ws 'B31' 'Filecoin'
Set-TextValue $ws 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D31' '4.530'
Set-TextValue $ws 'E31' '  +0.31%  '

# Row 32
Set-TextValue $ws 'B32' 'PancakeSwap'
Set-TextValue $ws 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D32' '1.528'
Set-TextValue $ws 'E32' '  -1.40%  '

# Row 33
Set-TextValue $ws 'D33' '4.258'
Set-TextValue $ws 'E33' '  +3.94%  '

# Row 34
Set-TextValue $ws 'D34' '0.05304'
Set-TextValue $ws 'E34' '  -3.13%  '

# Row 35
Set-TextValue $ws 'D35' '1.228'
Set-TextValue $ws 'E35' '  -1.35%  '

# Row 36
Set-TextValue $ws 'D36' '0.7476'
Set-TextValue $ws 'E36' '  -0.02%  '

# Row 37
Set-TextValue $ws 'E37' '  +0.17%  '

# Row 38
Set-TextValue $ws 'D38' '2.692'
Set-TextValue $ws 'E38' '  -0.03%  '

# Row 39
Set-TextValue $ws 'D39' '0.01941'
Set-TextValue $ws 'E39' '  +0.49%  '

# Row 40
Set-TextValue $ws 'D40' '2.753'
Set-TextValue $ws 'E40' '  -1.26%  '

# Row 41
Set-TextValue $ws 'B41' 'TheSandbox'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D41' '0.4456'
Set-TextValue $ws 'E41' '  -0.73%  '

# Row 42
Set-TextValue $ws 'B42' 'Maker'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D42' '1.097.23'
Set-TextValue $ws 'E42' '  +0.07%  '

# Row 43
Set-TextValue $ws 'D43' '6.078'
Set-TextValue $ws 'E43' '  +0.64%  '

# Row 44
Set-TextValue $ws 'D44' '71.94'
Set-TextValue $ws 'E44' '  -3.06%  '

# Row 45
Set-TextValue $ws 'D45' '0.8700'
Set-TextValue $ws 'E45' '  +2.26%  '

# Row 46
Set-TextValue $ws 'D46' '1.000'
Set-TextValue $ws 'E46' '  -0.01%  '

# Row 47
Set-TextValue $ws 'B47' 'Aptos'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D47' '7.726'
Set-TextValue $ws 'E47' '  +1.88%  '

# Row 48
Set-TextValue $ws 'B48' 'Quant'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D48' '102.49'
Set-TextValue $ws 'E48' '  +0.07%  '

# Row 49
Set-TextValue $ws 'B49' 'SynthetixNetwork'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue $ws 'D49' '3.051'
Set-TextValue $ws 'E49' '  +1.37%  '

# Row 50
Set-TextValue $ws 'B50' 'RenderToken'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D50' '1.829'
Set-TextValue $ws 'E50' '  -3.08%  '

# Row 51
Set-TextValue $ws 'D51' '2.024.78'
Set-TextValue $ws 'E51' '  -0.97%  '
